$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B","C","D","E","F","I","J","K","L","M")

$rowValues = @{
    2 = @(1.02, 1.049779759218475, 1.046719066332263, 1.055985252378272, 1.064647919388417, 1.037891377006613, 1.05481639696388, 1.049483715523721, 1.058724206674619, 1.067363290219442)
    3 = @(1.02, 1.051422459474209, 1.047919464890359, 1.057491640126349, 1.066341427089232, 1.038281166959796, 1.056105125808099, 1.050494804132188, 1.060042379528144, 1.068869846434747)
    4 = @(1.02, 1.05248305601586, 1.048693893065422, 1.058464461607182, 1.067435532105449, 1.038530648265959, 1.056936347906305, 1.051146174450458, 1.060892912643392, 1.069842515340153)
    5 = @(1.02, 1.052928382075462, 1.049018917043822, 1.058872987815109, 1.067895096139295, 1.03863487767489, 1.057285163016613, 1.051419330001378, 1.061249908183965, 1.070250915126878)
    6 = @(1.02, 1.053003122437076, 1.049073458217405, 1.058941555179441, 1.067972235961908, 1.038652340077506, 1.05734369388453, 1.051465154314116, 1.061309816226209, 1.070319457622475)
    7 = @(1.02, 1.052489008624616, 1.048698238185375, 1.058469922100866, 1.067441674369938, 1.038532043543998, 1.056941011256555, 1.051149827032658, 1.060897685056856, 1.069847974388837)
    8 = @(1.02, 1.050335409035238, 1.047125229401511, 1.056494745022864, 1.065220608550963, 1.038023677354945, 1.055252487491488, 1.049826016437621, 1.059170194438135, 1.067872894169933)
    9 = @(1.02, 1.046522046270962, 1.04433536220758, 1.052999143638482, 1.061293227389028, 1.037106751968268, 1.052256228533952, 1.047470986582579, 1.056107247268729, 1.06437544449139)
    10 = @(1.02, 1.04396667154221, 1.042462862256719, 1.050657977241498, 1.058665140362748, 1.036481066908011, 1.050244127970698, 1.045885527599277, 1.054052009273655, 1.062031666704186)
    11 = @(1.02, 1.042856887425301, 1.041648957671883, 1.049641533644228, 1.05752466011835, 1.036206676841025, 1.049369277515472, 1.045195245158143, 1.053158791693065, 1.061013760375718)
    12 = @(1.02, 1.042444155707767, 1.041346162354366, 1.049263563034532, 1.057100645555484, 1.036104231654054, 1.049043767548334, 1.04493826872574, 1.052826506058558, 1.060635196146535)
    13 = @(1.02, 1.042532711301673, 1.04141113461779, 1.049344658133246, 1.057191615961115, 1.03612623029426, 1.04911361567897, 1.044993417193804, 1.0528978054815, 1.06071642086142)
    14 = @(1.02, 1.042822781348342, 1.041623938258286, 1.04961029907727, 1.057489618965766, 1.036198219413001, 1.04934238207596, 1.045174015182269, 1.053131335222799, 1.060982477768453)
    15 = @(1.02, 1.043001435271375, 1.041754990396264, 1.049773913450694, 1.057673176439891, 1.036242504662504, 1.049483259184321, 1.045285211136155, 1.053275153369447, 1.061146341753088)
    16 = @(1.02, 1.044040253733997, 1.04251681233129, 1.050725377275118, 1.058740776421532, 1.036499203941616, 1.050302112175512, 1.045931259241769, 1.05411121899995, 1.062099156784269)
    17 = @(1.02, 1.044690986306284, 1.04299384674108, 1.051321473201501, 1.059409774450106, 1.036659294297539, 1.05081478606875, 1.04633549319946, 1.054634773281223, 1.062696011111117)
    18 = @(1.02, 1.045070231135891, 1.043271794388823, 1.051668905883737, 1.059799749278852, 1.036752338313102, 1.051113473710274, 1.046570913038085, 1.054939837348564, 1.063043854152571)
    19 = @(1.02, 1.045199490570631, 1.04336651700809, 1.051787327715473, 1.059932680251583, 1.036784007389951, 1.051215260059891, 1.046651123810353, 1.055043803000409, 1.063162410490579)
    20 = @(1.02, 1.044621201696357, 1.042942696400481, 1.051257544729012, 1.059338022257549, 1.036642152691697, 1.050759816918101, 1.046292160329785, 1.054578633657541, 1.062632004582673)
    21 = @(1.02, 1.042737377100782, 1.041561286074653, 1.049532086082275, 1.057401875368557, 1.036177034923129, 1.049275031388861, 1.045120849513275, 1.053062580561337, 1.060904143674523)
    22 = @(1.02, 1.041549992742403, 1.040689986220392, 1.048444796090101, 1.056182283052372, 1.035881560302275, 1.048338291107298, 1.044381069942321, 1.052106451256122, 1.059815052077418)
    23 = @(1.02, 1.042179731971778, 1.041152142752546, 1.04902142307961, 1.056829030818245, 1.036038486182639, 1.048835181583011, 1.04477355965008, 1.052613594655006, 1.06039266187603)
    24 = @(1.02, 1.044652735335577, 1.042965809959945, 1.05128643205733, 1.059370444727667, 1.036649899277541, 1.050784656181427, 1.046311741710012, 1.054604001714299, 1.062660927279984)
    25 = @(1.02, 1.04751015266267, 1.045058793808688, 1.053904691922899, 1.062310230743898, 1.037346322542319, 1.053033363803286, 1.048082506830726, 1.056901389603565, 1.065281712089731)
}

foreach ($row in $rowValues.Keys) {
    $values = $rowValues[$row]
    for ($i = 0; $i -lt $columns.Count; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $values[$i]
    }
}
